$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their original text representation (avoid Excel
# auto-converting numeric-looking strings into numbers/dates),
# matching the source workbook which stores these as inline strings.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "96.846.39"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.67%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.685.65"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.06%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.28"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.73%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +3.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "654.78"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.73%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.70%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.10%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.684.80"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.19%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "44.05"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.92%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000300"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +12.41%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.74"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.373.78"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "96.699.28"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.64%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "9.03"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.94%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.675.03"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.81%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.97"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.75%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.66"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.508"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -4.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "521.41"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.80%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.64%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +4.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.89"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.206"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +25.38%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "101.04"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.85%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "13.31"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +3.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "12.40"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.00"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.40%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.01%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.85"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.33%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.34%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "32.14"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.93%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "643.87"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +4.35%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.590"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.82%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.79"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.44%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +11.79%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "ImmutableX"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.04"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +6.26%  "
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.69"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -4.95%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.950"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.462"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +9.01%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.72%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.27"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.52"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.17%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.62%  "
